$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1669.2858
$ws.Range("I19").Value = 1500.5
$ws.Range("J19").Value = 1736.8
$ws.Range("K19").Value = 1500.5
$ws.Range("L19").Value = 1736.8
$ws.Range("M19").Value = -1325.5
$ws.Range("N19").Value = -2086.8
$ws.Range("H28").Value = 9122.857
$ws.Range("I28").Value = 496.1111
$ws.Range("K28").Value = 496.1111
$ws.Range("M28").Value = -11.11110000000002
$ws.Range("H32").Value = 1304.5
$ws.Range("I32").Value = 1966.6666
$ws.Range("J32").Value = 1123.909
$ws.Range("K32").Value = 1966.6666
$ws.Range("L32").Value = 1123.909
$ws.Range("M32").Value = -1640.6666
$ws.Range("N32").Value = -1775.909
$ws.Range("H40").Value = 1231.6666
$ws.Range("I40").Value = 1097.5
$ws.Range("K40").Value = 1097.5
$ws.Range("M40").Value = -922.5
$ws.Range("H43").Value = 1956.0834
$ws.Range("J43").Value = 2434.125
$ws.Range("L43").Value = 2434.125
$ws.Range("N43").Value = -2572.125
$ws.Range("I51").Value = 2000
$ws.Range("K51").Value = 2000
$ws.Range("M51").Value = -1516
$ws.Range("H53").Value = 55555770
$ws.Range("I53").Value = 238.27272
$ws.Range("J53").Value = 142857310
$ws.Range("K53").Value = 238.27272
$ws.Range("L53").Value = 142857310
$ws.Range("M53").Value = 398.72728
$ws.Range("N53").Value = -142858584
$ws.Range("H64").Value = 3297.0715
$ws.Range("I64").Value = 3297.5
$ws.Range("J64").Value = 3294.5
$ws.Range("K64").Value = 3297.5
$ws.Range("L64").Value = 3294.5
$ws.Range("M64").Value = -3049.5
$ws.Range("N64").Value = -3790.5
$ws.Range("H67").Value = 3297.0715
$ws.Range("I67").Value = 3297.5
$ws.Range("J67").Value = 3294.5
$ws.Range("K67").Value = 3297.5
$ws.Range("L67").Value = 3294.5
$ws.Range("M67").Value = -2439.5
$ws.Range("N67").Value = -5010.5
$ws.Range("H70").Value = 1624541.8
$ws.Range("I70").Value = 2526378.2
$ws.Range("J70").Value = 1236
$ws.Range("K70").Value = 7579134.600000001
$ws.Range("L70").Value = 3708
$ws.Range("M70").Value = -7578864.600000001
$ws.Range("N70").Value = -4248
$ws.Range("H73").Value = 1624541.8
$ws.Range("I73").Value = 2526378.2
$ws.Range("J73").Value = 1236
$ws.Range("K73").Value = 7579134.600000001
$ws.Range("L73").Value = 3708
$ws.Range("M73").Value = -7578198.600000001
$ws.Range("N73").Value = -5580
$ws.Range("H74").Value = 5453.55
$ws.Range("I74").Value = 5485.7144
$ws.Range("J74").Value = 5378.5
$ws.Range("K74").Value = 5485.7144
$ws.Range("L74").Value = 5378.5
$ws.Range("M74").Value = -4549.7144
$ws.Range("N74").Value = -7250.5
$ws.Range("H76").Value = 3068.0256
$ws.Range("I76").Value = 3013.7273
$ws.Range("K76").Value = 3013.7273
$ws.Range("M76").Value = -2698.7273
$ws.Range("H77").Value = 5453.55
$ws.Range("I77").Value = 5485.7144
$ws.Range("J77").Value = 5378.5
$ws.Range("K77").Value = 27428.572
$ws.Range("L77").Value = 26892.5
$ws.Range("M77").Value = -22748.572
$ws.Range("N77").Value = -36252.5
$ws.Range("H79").Value = 3068.0256
$ws.Range("I79").Value = 3013.7273
$ws.Range("K79").Value = 3013.7273
$ws.Range("M79").Value = -1921.7273
$ws.Range("H80").Value = 2789.9167
$ws.Range("I80").Value = 3725.5715
$ws.Range("J80").Value = 1480
$ws.Range("K80").Value = 11176.7145
$ws.Range("L80").Value = 4440
$ws.Range("M80").Value = -10178.7145
$ws.Range("N80").Value = -6436
$ws.Range("H83").Value = 2789.9167
$ws.Range("I83").Value = 3725.5715
$ws.Range("J83").Value = 1480
$ws.Range("K83").Value = 33530.1435
$ws.Range("L83").Value = 13320
$ws.Range("M83").Value = -28538.1435
$ws.Range("N83").Value = -23304
$ws.Range("H86").Value = 210021550
$ws.Range("I86").Value = 262526130
$ws.Range("J86").Value = 3250
$ws.Range("K86").Value = 262526130
$ws.Range("L86").Value = 3250
$ws.Range("M86").Value = -262525007
$ws.Range("N86").Value = -5496
$ws.Range("H88").Value = 77473.5
$ws.Range("I88").Value = 1690
$ws.Range("J88").Value = 102734.664
$ws.Range("K88").Value = 1690
$ws.Range("L88").Value = 102734.664
$ws.Range("M88").Value = -1284
$ws.Range("N88").Value = -103546.664
$ws.Range("H89").Value = 210021550
$ws.Range("I89").Value = 262526130
$ws.Range("J89").Value = 3250
$ws.Range("K89").Value = 1312630650
$ws.Range("L89").Value = 16250
$ws.Range("M89").Value = -1312625034
$ws.Range("N89").Value = -27482
$ws.Range("H91").Value = 77473.5
$ws.Range("I91").Value = 1690
$ws.Range("J91").Value = 102734.664
$ws.Range("K91").Value = 1690
$ws.Range("L91").Value = 102734.664
$ws.Range("M91").Value = -286
$ws.Range("N91").Value = -105542.664
$ws.Range("H100").Value = 2698.6316
$ws.Range("I100").Value = 2638.5
$ws.Range("J100").Value = 2765.4443
$ws.Range("K100").Value = 2638.5
$ws.Range("L100").Value = 2765.4443
$ws.Range("M100").Value = -2097.5
$ws.Range("N100").Value = -3847.4443
$ws.Range("H107").Value = 1292.8334
$ws.Range("I107").Value = 1427.9231
$ws.Range("K107").Value = 1427.9231
$ws.Range("M107").Value = 492.0769
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = $null
$ws.Range("N119").Value = 0
$ws.Range("H137").Value = 2401.5
$ws.Range("I137").Value = 2501.6072
$ws.Range("J137").Value = 1000
$ws.Range("K137").Value = 7504.821599999999
$ws.Range("L137").Value = 3000
$ws.Range("M137").Value = -4954.821599999999
$ws.Range("N137").Value = -8100
$ws.Range("H138").Value = 219536.98
$ws.Range("I138").Value = 3115.2307
$ws.Range("J138").Value = 295577.06
$ws.Range("K138").Value = 9345.6921
$ws.Range("L138").Value = 886731.1799999999
$ws.Range("M138").Value = -4205.6921
$ws.Range("N138").Value = -897011.1799999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1903.3334
$ws.Range("I102").Value = 1903.3334
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1903.3334
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = $null
$ws.Range("N102").Value = -281.3334
$ws.Range("H123").Value = 35762
$ws.Range("J123").Value = 35762
$ws.Range("L123").Value = 35762
$ws.Range("N123").Value = -45562

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1900
$ws.Range("I99").Value = 1900
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1900
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = $null
$ws.Range("N99").Value = -402
$ws.Range("H107").Value = 1721.3636
$ws.Range("I107").Value = 941.8
$ws.Range("K107").Value = 941.8
$ws.Range("M107").Value = 978.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2867.889
$ws.Range("I132").Value = 1100
$ws.Range("J132").Value = 3373
$ws.Range("K132").Value = 3300
$ws.Range("L132").Value = 10119
$ws.Range("M132").Value = -770
$ws.Range("N132").Value = -15179
$ws.Range("H133").Value = 50313
$ws.Range("J133").Value = 50313
$ws.Range("L133").Value = 50313
$ws.Range("N133").Value = -55373

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 1400
$ws.Range("I63").Value = 880
$ws.Range("K63").Value = 2640
$ws.Range("M63").Value = -1891
$ws.Range("H66").Value = 1400
$ws.Range("I66").Value = 880
$ws.Range("K66").Value = 7920
$ws.Range("M66").Value = -4176
$ws.Range("H122").Value = 8247.154
$ws.Range("I122").Value = 423.8889
$ws.Range("J122").Value = 25849.5
$ws.Range("K122").Value = 3815.0001
$ws.Range("L122").Value = 232645.5
$ws.Range("M122").Value = -1365.0001
$ws.Range("N122").Value = -237545.5
$ws.Range("H132").Value = 3089.2683
$ws.Range("I132").Value = 2069
$ws.Range("K132").Value = 18621
$ws.Range("M132").Value = -16091
$ws.Range("H137").Value = 11913954
$ws.Range("I137").Value = 13899121
$ws.Range("K137").Value = 41697363
$ws.Range("M137").Value = -41692263

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 55000
$ws.Range("J133").Value = 55000
$ws.Range("L133").Value = 55000
$ws.Range("N133").Value = -65120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 39003.5
$ws.Range("I24").Value = 8000
$ws.Range("K24").Value = 8000
$ws.Range("M24").Value = -7657
$ws.Range("H100").Value = 2543
$ws.Range("I100").Value = 2449.8333
$ws.Range("J100").Value = 2822.5
$ws.Range("K100").Value = 2449.8333
$ws.Range("L100").Value = 2822.5
$ws.Range("M100").Value = -1908.8333
$ws.Range("N100").Value = -3904.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 841.5
$ws.Range("I100").Value = 409.8
$ws.Range("K100").Value = 819.6
$ws.Range("M100").Value = -278.6
$ws.Range("H123").Value = 30856.25
$ws.Range("J123").Value = 30856.25
$ws.Range("L123").Value = 30856.25
$ws.Range("N123").Value = -40656.25
$ws.Range("H126").Value = 1107.9667
$ws.Range("I126").Value = 760.5263
$ws.Range("J126").Value = 1708.091
$ws.Range("K126").Value = 2281.5789
$ws.Range("L126").Value = 5124.272999999999
$ws.Range("M126").Value = 188.4211
$ws.Range("N126").Value = -10064.273
